$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Change 1: fix subject/verb agreement ("se encuentran" -> "se encuentra")
# ---------------------------------------------------------------
$d.Content.Find.Execute(
    "El equipo inspeccionado, identificado en el ítem II, ubicado en {{inspection_place}}, durante la inspección se puede apreciar que se encuentran en {{buen/mal}} estado de operación.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "El equipo inspeccionado, identificado en el ítem II, ubicado en {{inspection_place}}, durante la inspección se puede apreciar que se encuentra en {{buen/mal}} estado de operación.",
    2) | Out-Null

# ---------------------------------------------------------------
# Change 2: append an inspection-method note block after the closing
# signature table, starting on a new page.
# ---------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$tailRange = $lastPara.Range.Duplicate()
$tailRange.MoveEnd(1, -1) | Out-Null
$tailRange.InsertAfter([char]12) | Out-Null

# New paragraph: "MÉTODO DE INSPECCION: "
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$p1 = $d.Paragraphs.Last
$p1.Alignment = 1
$p1.Range.Font.Bold = $true
$p1.Range.Font.BoldBi = $true
$p1.Range.Font.Size = 10
$p1.Range.Font.SizeBi = 10
$p1.Range.InsertAfter("MÉTODO DE INSPECCION: ") | Out-Null

# New paragraph: "DO: Documental; VI: Visual; FU: Funcionamiento; DI: Dimensional."
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$p2 = $d.Paragraphs.Last
$p2.Alignment = 1
$p2.Range.Font.Bold = $true
$p2.Range.Font.BoldBi = $true
$p2.Range.Font.Size = 10
$p2.Range.Font.SizeBi = 10
$p2.Range.InsertAfter("DO: Documental; VI: Visual; FU: Funcionamiento; DI: Dimensional.") | Out-Null

# New (empty) trailing paragraph, keeps the same bold/centered formatting
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$p3 = $d.Paragraphs.Last
$p3.Alignment = 1
$p3.Range.Font.Bold = $true
$p3.Range.Font.BoldBi = $true
$p3.Range.Font.Size = 10
$p3.Range.Font.SizeBi = 10
